# Update gh-pages to output generated at 456a3b4
# Applies updated attendee/follower counts (column F) across the
# "展览" (sheet 1), "演出" (sheet 2) and "全部类型" (sheet 4) worksheets.

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item(1)
$sheet1Updates = @{
    5  = 209
    6  = 1137
    7  = 1083
    8  = 8412
    11 = 6975
    13 = 306
    14 = 5135
    15 = 5135
    18 = 5646
    19 = 5647
    20 = 1081
    22 = 360
    24 = 490
    28 = 9446
    30 = 1736
    31 = 1147
    33 = 1919
    35 = 80
    36 = 149
    37 = 1020
    38 = 1933
    39 = 249
    40 = 1246
    41 = 56
    42 = 4924
    44 = 1171
    45 = 539
    50 = 1289
}
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Cells.Item($row, 6).Value = $sheet1Updates[$row]
}

# Sheet 2: 演出
$ws2 = $wb.Worksheets.Item(2)
$sheet2Updates = @{
    10 = 37
}
foreach ($row in $sheet2Updates.Keys) {
    $ws2.Cells.Item($row, 6).Value = $sheet2Updates[$row]
}

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item(4)
$sheet4Updates = @{
    6  = 209
    8  = 1137
    9  = 1083
    10 = 8412
    13 = 6975
    15 = 306
    18 = 5135
    19 = 5135
    21 = 5647
    22 = 5647
    23 = 1081
    25 = 360
    26 = 490
    31 = 9446
    33 = 1736
    34 = 1148
    36 = 1919
    38 = 80
    39 = 1020
    40 = 1933
    41 = 249
    42 = 1246
    43 = 4925
    45 = 1171
    46 = 539
    51 = 1289
}
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Cells.Item($row, 6).Value = $sheet4Updates[$row]
}
